$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "roboticS1Prep" column (I) was stored as the literal string "No" in
# every data row. Convert it to a real boolean (FALSE) with a custom
# TRUE/FALSE display format, and move the active selection from the old
# H2:H32 block onto the now-relevant I2:I32 block.
$rng = $ws.Range("I2:I32")
$rng.Value = $false
$rng.NumberFormat = '"TRUE";"TRUE";"FALSE"'
$rng.Select()
